# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Ajo".
# The new data point is inserted at row 241, pushing the existing rows
# (old 241..314) down to 242..315 and extending the sheet's used range
# from A1:R314 to A1:R315.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 241..314 down by one row, carrying formatting (incl. the
# date-format style on column D) along with them.
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new observation.
$ws.Cells.Item(241, 1).Value  = 4
$ws.Cells.Item(241, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(241, 3).Value  = "Los Lagos"
$ws.Cells.Item(241, 4).Value  = 44809
$ws.Cells.Item(241, 5).Value  = 10
$ws.Cells.Item(241, 6).Value  = 100112003
$ws.Cells.Item(241, 7).Value  = "Ajo"
$ws.Cells.Item(241, 8).Value  = "Chino"
$ws.Cells.Item(241, 9).Value  = "Primera"
$ws.Cells.Item(241, 10).Value = 60
$ws.Cells.Item(241, 11).Value = 29000
$ws.Cells.Item(241, 12).Value = 29000
$ws.Cells.Item(241, 13).Value = 29000
$ws.Cells.Item(241, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(241, 15).Value = "China"
$ws.Cells.Item(241, 16).Value = 2900
$ws.Cells.Item(241, 17).Value = 10
$ws.Cells.Item(241, 18).Value = "Hortaliza"
